$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -538.9699999999993
$ws.Range("C2").Value = 94.75
$ws.Range("D2").Value = 954.2000000000002
$ws.Range("B3").Value = 363.3600000000001
$ws.Range("C3").Value = 156.21
$ws.Range("D3").Value = 1391.0300000000002
$ws.Range("B4").Value = 1361.1599999999999
$ws.Range("C4").Value = 161.07
$ws.Range("D4").Value = 1923.2400000000002
$ws.Range("B5").Value = 2186.06
$ws.Range("C5").Value = 159.87
$ws.Range("D5").Value = 2337.17
$ws.Range("B6").Value = 2931.579999999999
$ws.Range("C6").Value = 187.57000000000002
$ws.Range("D6").Value = 2719.7000000000003
$ws.Range("B7").Value = 3584.1899999999996
$ws.Range("C7").Value = 207.51000000000002
$ws.Range("D7").Value = 3040.68
$ws.Range("B8").Value = 4107.5
$ws.Range("C8").Value = 225.18
$ws.Range("D8").Value = 3280.7000000000007
$ws.Range("B9").Value = 4567.849999999999
$ws.Range("C9").Value = 239.12
$ws.Range("D9").Value = 3497.66
$ws.Range("B10").Value = 5037.59
$ws.Range("C10").Value = 246.79
$ws.Range("D10").Value = 3706.58
$ws.Range("B11").Value = 5462.8099999999995
$ws.Range("C11").Value = 270.65999999999997
$ws.Range("D11").Value = 3888.37
$ws.Range("B12").Value = 5836.67
$ws.Range("C12").Value = 283.36
$ws.Range("D12").Value = 4036.09
$ws.Range("B13").Value = 6170.52
$ws.Range("C13").Value = 287.77
$ws.Range("D13").Value = 4157.16
$ws.Range("B14").Value = 6479.64
$ws.Range("C14").Value = 303.05
$ws.Range("D14").Value = 4269.740000000001
$ws.Range("B15").Value = 6830.59
$ws.Range("C15").Value = 319.5
$ws.Range("D15").Value = 4374.25
$ws.Range("B16").Value = 7188.959999999999
$ws.Range("C16").Value = 328.65000000000003
$ws.Range("D16").Value = 4486.88
$ws.Range("B17").Value = 7574.22
$ws.Range("C17").Value = 338.41999999999996
$ws.Range("D17").Value = 4608.66
$ws.Range("B18").Value = 7898.35
$ws.Range("C18").Value = 346.40999999999997
$ws.Range("D18").Value = 4689.0
$ws.Range("B19").Value = 8030.16
$ws.Range("C19").Value = 358.23
$ws.Range("D19").Value = 4672.92
$ws.Range("B20").Value = 8182.66
$ws.Range("C20").Value = 379.33
$ws.Range("D20").Value = 4645.950000000001
$ws.Range("B21").Value = 8542.85
$ws.Range("C21").Value = 388.84999999999997
$ws.Range("D21").Value = 4723.32
$ws.Range("B22").Value = 8739.43
$ws.Range("C22").Value = 398.53999999999996
$ws.Range("D22").Value = 4737.43
$ws.Range("B23").Value = 9069.1
$ws.Range("C23").Value = 412.0
$ws.Range("D23").Value = 4786.5599999999995
$ws.Range("B24").Value = 9199.8
$ws.Range("C24").Value = 419.46
$ws.Range("D24").Value = 4733.09
$ws.Range("B25").Value = 9420.789999999999
$ws.Range("C25").Value = 449.3
$ws.Range("D25").Value = 4730.119999999999
$ws.Range("B26").Value = 9679.41
$ws.Range("C26").Value = 451.22
$ws.Range("D26").Value = 4698.35
$ws.Range("B27").Value = 9865.41
$ws.Range("C27").Value = 455.95
$ws.Range("D27").Value = 4668.22
$ws.Range("B28").Value = 9928.07
$ws.Range("C28").Value = 466.71
$ws.Range("D28").Value = 4577.849999999999
$ws.Range("B29").Value = 10125.560000000001
$ws.Range("C29").Value = 491.65999999999997
$ws.Range("D29").Value = 4527.709999999999
$ws.Range("B30").Value = 10259.949999999999
$ws.Range("C30").Value = 489.78999999999996
$ws.Range("D30").Value = 4449.08
$ws.Range("B31").Value = 10474.83
$ws.Range("C31").Value = 505.46
$ws.Range("D31").Value = 4395.66
$ws.Range("B32").Value = 10490.27
$ws.Range("C32").Value = 529.64
$ws.Range("D32").Value = 4246.07
$ws.Range("B33").Value = 10895.509999999998
$ws.Range("C33").Value = 537.28
$ws.Range("D33").Value = 4271.14
$ws.Range("B34").Value = 11075.050000000001
$ws.Range("C34").Value = 536.31
$ws.Range("D34").Value = 4173.95
$ws.Range("B35").Value = 11234.28
$ws.Range("C35").Value = 566.24
$ws.Range("D35").Value = 4078.3900000000003
$ws.Range("B36").Value = 11394.43
$ws.Range("C36").Value = 575.38
$ws.Range("D36").Value = 3970.3599999999997
$ws.Range("B37").Value = 11405.18
$ws.Range("C37").Value = 568.8000000000001
$ws.Range("D37").Value = 3826.9100000000003
$ws.Range("B38").Value = 11579.47
$ws.Range("C38").Value = 588.83
$ws.Range("D38").Value = 3716.06
$ws.Range("B39").Value = 11756.0
$ws.Range("C39").Value = 609.14
$ws.Range("D39").Value = 3609.26
$ws.Range("B40").Value = 11915.11
$ws.Range("C40").Value = 618.51
$ws.Range("D40").Value = 3500.57
$ws.Range("B41").Value = 12107.68
$ws.Range("C41").Value = 628.5600000000001
$ws.Range("D41").Value = 3356.3999999999996
$ws.Range("B42").Value = 12069.259999999998
$ws.Range("C42").Value = 645.9399999999999
$ws.Range("D42").Value = 3156.66
$ws.Range("B43").Value = 12150.75
$ws.Range("C43").Value = 657.6999999999999
$ws.Range("D43").Value = 3018.5
$ws.Range("B44").Value = 12425.699999999999
$ws.Range("C44").Value = 665.84
$ws.Range("D44").Value = 2897.75
$ws.Range("B45").Value = 12534.359999999999
$ws.Range("C45").Value = 678.37
$ws.Range("D45").Value = 2742.99
$ws.Range("B46").Value = 12733.689999999999
$ws.Range("C46").Value = 699.83
$ws.Range("D46").Value = 2598.21
$ws.Range("B47").Value = 12675.66
$ws.Range("C47").Value = 705.6100000000001
$ws.Range("D47").Value = 2404.95
$ws.Range("B48").Value = 12973.710000000001
$ws.Range("C48").Value = 725.72
$ws.Range("D48").Value = 2273.68
$ws.Range("B49").Value = 13002.61
$ws.Range("C49").Value = 718.37
$ws.Range("D49").Value = 2068.99
$ws.Range("B50").Value = 12942.98
$ws.Range("C50").Value = 734.73
$ws.Range("D50").Value = 1881.04
$ws.Range("B51").Value = 12993.060000000001
$ws.Range("C51").Value = 738.97
$ws.Range("D51").Value = 1689.88
$ws.Range("B52").Value = 13107.260000000002
$ws.Range("C52").Value = 740.86
$ws.Range("D52").Value = 1525.5399999999997
$ws.Range("B53").Value = 13001.300000000001
$ws.Range("C53").Value = 748.64
$ws.Range("D53").Value = 1282.6999999999998
$ws.Range("B54").Value = 13113.61
$ws.Range("C54").Value = 750.36
$ws.Range("D54").Value = 1109.9299999999998
$ws.Range("B55").Value = 13317.5
$ws.Range("C55").Value = 761.04
$ws.Range("D55").Value = 933.4800000000001
$ws.Range("B56").Value = 13432.73
$ws.Range("C56").Value = 784.22
$ws.Range("D56").Value = 764.4600000000003
$ws.Range("B57").Value = 13218.76
$ws.Range("C57").Value = 771.4100000000001
$ws.Range("D57").Value = 541.8700000000001
$ws.Range("B58").Value = 13453.18
$ws.Range("C58").Value = 793.2299999999999
$ws.Range("D58").Value = 324.0299999999999

$ws.Range("O15").Select()

